# Update the "calcTime" (column E) measurements for each of the four
# method blocks on the metrics sheet. These values reflect re-running the
# benchmark after reworking coordsToU (start/end point visualization +
# togglable figure generation) and aligning multiObstacleReset's call
# signature with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value  = 0.1643925
$ws.Range("E4").Value  = 0.015528999999999999
$ws.Range("E5").Value  = 0.023719299999999999

$ws.Range("E11").Value = 0.012711099999999999
$ws.Range("E12").Value = 0.0225343
$ws.Range("E13").Value = 0.031447599999999999

$ws.Range("E19").Value = 0.0103039
$ws.Range("E20").Value = 0.011632200000000001
$ws.Range("E21").Value = 0.025541600000000001

$ws.Range("E27").Value = 0.0083443000000000007
$ws.Range("E28").Value = 0.010261599999999999
$ws.Range("E29").Value = 0.026808499999999999
